# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AA1) so the
# new headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the season record for every data row (2-45).
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}
